$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 're : def retirement and employee scheme - abc good morning further to your email please see the following : 10 test road test street aa12 123 1st january 1400 zz123123z regards test ifa from : test pensions pensions com sent : 01 december 2020 08 : 04 to : test ifa fincorp com subject : def retirement and employee scheme abc dear test ifa def retirement and employee scheme thank you for your email regarding the def retirement and employee scheme to enable us to answer your enquiry please could you provide the member full address date of birth and national insurance number for data verification purposes for added security please also confirm your full company name and the address stated on the provided letter of authority should you have any further queries please do not hesitate to contact the customer support centre between 8 30am and 5 30pm monday to friday on our helpline + 91 012345678 alternatively please email us at test pensions pensions com yours sincerely test pensions from : test ifa fincorp com sent : 27 november 2020 14 : 55 to : test pensions pensions com subject : - tst12345678 good morning i uploaded an loa and information request form to your website on 29 10 20 please could you provide me with an update kind regards test ifa'
$ws.Cells.Item(2, 2).Value = 'DEF Retirement and Employee Scheme - ABC.msg'

$ws.Cells.Item(3, 1).Value = 'letter of authority - client abcd good afternoon please find attached a letter of authority and information request for our client abcd please send all correspondence to this email where possible kind regards test ifa'
$ws.Cells.Item(3, 2).Value = 'Letter of Authority - client ABCD.msg'

$ws.Cells.Item(4, 1).Value = 'fw : cetv request mr xyz following my email dated 28th october mr xyz has now left employment and is no longer an active member of the pension scheme please calculate and send his guaranteed transfer value cetv as soon as possible regards test ifa from : test ifa fincorp com sent : 28 october 2020 13 : 45 to : test pensions pensions com subject : cetv request mr xyz importance : high mr xyz 10 test road test street aa12 123 works number : tst12345678 dob : 01 01 1400 ni no : zz123123z i hold authority for mr xyz and request a cetv be sent to me for him as soon as possible thank you test ifa'
$ws.Cells.Item(4, 2).Value = 'CETV Request_Mr XYZ .msg'

$ws.Cells.Item(5, 1).Value = 'fwd : early quotation get outlook for android from : abc xyz sent : tuesday august 27 2019 10 : 55 : 53 am to : test ifa fincorp com subject : early quotation please can you send me some early retirement quotation forms ref tst12345678 national insurance number zz123123z get outlook for android'
$ws.Cells.Item(5, 2).Value = 'FW Early quotation.msg'

$ws.Cells.Item(6, 1).Value = 'fw : information request - letter of authority attached ref : tst12345678 hello team please see below query and reply accordingly to the member ifa many thanks test pensions from : test ifa fincorp com sent : 01 december 2020 08 : 35 to : test pensions pensions com subject : information request - letter of authority attached ref : tst12345678 dear representative of xx please see attached herewith an authorisation letter from with the following personal information : customer personal information : name : abc current address : 10 test road test street aa12 123 date of birth : 01 01 1400 gender : male national insurance number : zz123123z pension type : workplace workplace pension information : employer name : test employer pension provider : xx pension reference number : not known estimated valuation : not known additional information provided : please see information request attachment please note any reference number s or scheme name s listed above are those known to the customer there may however be additional pensions belonging to the customer therefore please ensure you provide information on all of test pensions with xx not just those referenced above with respect to the customer above please provide information on the following : latest transfer valuation current status of the plan : in force paid up or in payment are there any exit charges associated with the specific funds or with the plan overall any guaranteed or non - standard benefits including guaranteed salaries guaranteed annuity rates guaranteed maturity amount guaranteed minimum pension and any other benefits and the conditions attached to them please confirm whether partial transfers are possible from this plan if transfer forms are required to action a transfer out of the scheme please include these in your response test ifa is a company with the following details : company name : test ifa email : test ifa fincorp com address : 10 test road test street aa12 123 telephone : + 91 012345678 please contact test ifa if you require further information sincerely test ifa'
$ws.Cells.Item(6, 2).Value = 'Information Request - Letter of Authority attached (Ref TST12345678).msg'

$ws.Cells.Item(7, 1).Value = ' pre - identified user support request reply email address : test pensions pensions com forename : abc xyz surname : def ni number : zz123123z dob : 01 jan 1400 home postcode : aa12 123 employer : fin corp ltd subject : question re my test pension and the lifetimes allowance content : i am 55 in 1400 and i am now looking to start considering my pension retirement options i have a sipp with over 100 000 in it so looking ahead i am conscious that over the next few years i may come close to the lifetime allowance when i also consider the value of my test pension my understanding is that to work out the amount of my test pension within the lifetime allowance i multiply the annual pension payment by 20 is it possible to tell me what my pension is for lifetimes allowance calculation at 55 amd 60 much appreciated'
$ws.Cells.Item(7, 2).Value = 'FW PRE-IDENTIFIED USER SUPPORT REQUEST 2.msg'

$ws.Cells.Item(8, 1).Value = 'fw : tst12345678 thanks test pensions i have copied all the admin details on our data so someone from the team can look into this though looking at the file there does seem to be a few emails on the file to ifa so may just be a lack of clarity on what the hold up for mr abc is please can someone contact the client asap to nip this complaint in the bud thanks test fund - - - - - original message - - - - - from : test pensions pensions com sent : 01 december 2020 10 : 11 to : testfund pensions com subject : fw : tst12345678 hi test fund i am not sure if stu is around but please can you pass this on to the manager of mr abc scheme please i have issued an initial holding response as per the attached but this will need to be escalated accordingly thanks test pensions - - - - - original message - - - - - from : test pensions pensions com sent : 01 december 2020 10 : 04 to : abc xyz testmail com subject : re : tst12345678 good morning mr abc firstly please accept my apologies on behalf of xyz for the problems that you have been experiencing we take matters like this very seriously and i will do everything i can to assist you your email has come through to our specialist sipp and ssas team mailbox and we are unfortunately not the area that look after your pension scheme we also do not have access to the main pensions administration area records to enable me to search for your scheme and trace your policy however i have referred your email across to our complaints team and also to a senior manager in our main pensions administration area to enable them to trace your scheme and make contact with you regarding this they should be in touch with you shortly in the meantime i can also confirm that there is a customer service centre for our main pensions administration area which you can contact if you need to in the future by telephone on + 91 012345678 or by email to test pensions pensions com to help them to trace your scheme in the future you may need to quote the name of the scheme your policy reference number and your name and national insurance number i hope this is of assistance to you for the future and please rest assured that someone will be in touch with you shortly regarding your scheme kind regards test pensions - - - - - original message - - - - - from : abc xyz testmail com sent : 30 november 2020 14 : 30 to : test pensions pensions com subject : tst12345678 hi i ve sent the below message to the contact us page on your website someone needs to contact me or reply to test ifa with regards to transferring this pension to them you do not respond to me you do not respond to them nobody answers your phone this request was first made in september and will be sent to the ombudsman this week unless you actually reply to someone i will try here as well can someone respond to this before i m left no option but to go to the ombudsman from my point of view you have 20000 of my money you refuse to handover i consider that a form of theft i know someone looks at this as i ve been sent three useless letters in the post that achieve nothing you need to reply to test ifa when they contact you i thought def had the worst customer service ever but at least they eventually respond this is just corporate stonewalling have you tried using your service it s pathetic thanks abc'
$ws.Cells.Item(8, 2).Value = 'FW TST12345678.msg'

$ws.Cells.Item(9, 1).Value = 'fw : abcd query hello team please see below query and reply accordingly to the member ifa many thanks test pensions from : test ifa fincorp com sent : 30 november 2020 17 : 33 to : test pensions pensions com subject : abcd query hi please find attached a letter of authority and scheme information request letter that was previously sent to xyz in september 2020 can you please confirm receipt and arrange for this information to be provided to me as a matter of urgency we have received other cetv s for this client which are administered by yourselves which are due to expire at end of december 2020 therefore we would appreciate your assistance to expedite this request we have unfortunately spent some time trying to trace the administrators of this scheme and i understand the benefits were invested through xyz your assistance in this matter would be gratefully appreciated password to follow kind regards test ifa'
$ws.Cells.Item(9, 2).Value = 'ABCD QUERY.msg'

$ws.Cells.Item(10, 1).Value = 're : change of marital status - abc - zz123123z abc from : test jobcreation pensions com sent : 15 october 2020 14 : 09 to : test pensions pensions com subject : fw : change of marital status - abc - zz123123z hi test pensions please provide the email address as it is missing on below email regards test creation from : abc xyz testmail com sent : tuesday october 13 2020 2 : 47 : 37 to : test pensions pensions com subject : change of marital status - abc - zz123123z hi i would like to notify you that i recently married and would ask that you amend your records to ensure my husband xyz is now listed as the beneficiary to my pension i have attached a copy of my wedding certificate my address is : 10 test road test street aa12 123 please let me know if you require any further information many thanks and regards abc'
$ws.Cells.Item(10, 2).Value = 'Change of marital status - ABC - ZZ123123Z.msg'

$ws.Cells.Item(11, 1).Value = 'immediate retirement quote ref : tst12345678 test pension scheme mr abc xyz 10 test road test street aa12 123 ni zz123123z dear abc thank you for your recent correspondence i would like to request an immediate retirement quote for mr abc and retirement pack with his options on taking his pension from the scheme if you could also please send me a copy of this information if you require any further information from myself please let me know kind regards abc xyz'
$ws.Cells.Item(11, 2).Value = 'FW Immediate Retirement Quote.msg'

$ws.Cells.Item(12, 1).Value = 'change of address hi my name is abc my birthday is 01 01 1400 around a year ago i moved house my new address is as below mr abc 10 test road test street aa12 123 my previous address was 103 test road test street ua12 123 i have tried to register online but it didn t seem to work properly regards abc '
$ws.Cells.Item(12, 2).Value = 'Change of address_ABC.msg'

$ws.Cells.Item(13, 1).Value = 'abc - cetv request dear sirs re : test group pension plan member : abc ni number : zz123123z i refer to my request of 12 october 2020 for a new cetv in respect of the above member and your letter of 2 november confirming it was referred to the scheme actuary for review i should be grateful if you would provide me with an update on this and when i should expect receipt yours faithfully test fund'
$ws.Cells.Item(13, 2).Value = 'ABC - CETV Request.msg'

$ws.Cells.Item(14, 1).Value = 'change of name dear sir abc pension reference tst12345678 zz123123z i refer to a mail forwarded to your payroll department on 13th october on reflection i think this should have been sent to you i would like to confirm my change of name i married on 4th august 2020 and my new surname is abc i have advised hmrc and my details have been updated with them similarly my bank have been informed and my accounts have been amended from def to abc but the account number and sort code remain the same i attach a photocopy of my marriage certificate for your records should you require any further information please do not hesitate to contact me my phone numbers and email address remain the same yours faithfully mrs abc'
$ws.Cells.Item(14, 2).Value = 'Change of name.msg'

$ws.Cells.Item(15, 1).Value = 'pension hello my name is abc xyz i am 56 years old and would like to take a portion of my gkn pension but not sure how to go about it would be thankful if someone could contact me on 123456789 many thanks sent from my huawei phone'
$ws.Cells.Item(15, 2).Value = 'FW Pension .msg'

$ws.Cells.Item(16, 1).Value = 'fw : change of address for test pension - - - - - original message - - - - - from : abc xyz testmail com sent : 16 october 2020 12 : 59 to : test pensions pensions com subject : change of address for test pension dear sir ma am please note a change of address for a recipient of a test pension pension ref tst12345678 mrs abc of 10 test road test street aa12 123 has now moved to 10 test road test street aa12 123 banking details remain the same please update your records accordingly many thanks abc + 91 012345678'
$ws.Cells.Item(16, 2).Value = 'Change of address for TEST pension.msg'

$ws.Cells.Item(17, 1).Value = 'pension value dear sir madam i would like to know how i can review the total value of my pension and its transfer value i can access fin corp but can only find my annual entitlement at retirement age regards abc xyz supply chain me direct : + 91 012345678 mobile : + 91 012345678 email : abc xyz testmail com'
$ws.Cells.Item(17, 2).Value = 'FW Pension value.msg'

$ws.Cells.Item(18, 1).Value = 'pension abc xyz 01 01 1400 zz123123z hi i got a quote this morning to take my pension and want to know was my care allowance included in the quote and also i would like a quote to take my whole pension as one lump sum please can you provide this for me as i can t do this on line thanks in advance abc xyz'
$ws.Cells.Item(18, 2).Value = 'FW Pension.msg'

$ws.Cells.Item(19, 1).Value = 'policy ref : 13245 dear sirs re : policy ref 132456 mrs abc xyz 10 test road test street aa12 123 date of birth 01 01 1400 ni : zz123123z my employment at yorkshire water ceased on 01 01 1400 and i am writing to request early retirement quotation estimation calculations from you for the following options : 01 01 1400 with lump sum and without lump sum options - age 58 01 01 1400 with lump sum and without lump sum options - age 60 i look forward to hearing from you soon yours faithfully mrs abc xyz'
$ws.Cells.Item(19, 2).Value = 'FW Policy Ref 132456.msg'

$ws.Cells.Item(20, 1).Value = 'change of address details hi i am in the test pension scheme and we have moved house in the last few days i have updated my address details with hr here at xxx but wondered how i can do this on my test pension account or will this update automatically over time now that hr are aware regards abc'
$ws.Cells.Item(20, 2).Value = 'Change of address details.msg'

$ws.Cells.Item(21, 1).Value = ' pre - identified user support request reply email address : abc xyz testmail com forename : abc xyz surname : def ni number : zz123123z dob : 01 jan 1400 home postcode : aa12 123 employer : fin corp ltd subject : leaving fin corp ltd content : hi i will be leaving fin corp ltd in february 1400 but not taking any pension at that time my nrd is not until 1500 am i ok to simply leave the arrangements as they stand at the point when i leave or do i need to make any addtional arrangements thanks abc'
$ws.Cells.Item(21, 2).Value = 'FW PRE-IDENTIFIED USER SUPPORT REQUEST.msg'

$ws.Cells.Item(22, 1).Value = 'lmno request for a guaranteed cetv pack to test pensions regarding your existing client lmno please find attached the new letter of request for a guaranteed cetv pack for lmno our letter of authority is already on your files and the client has now opted out of the pension scheme please send all the relevant information by post and send a copy by return email please contact me if you need any further information to be able to send the guaranteed cetv pack to us for lmno kind regards test ifa 10 test road test street aa12 123 + 91 012345678'
$ws.Cells.Item(22, 2).Value = 'LMNO request for a guaranteed CETV pack.msg'

$ws.Cells.Item(23, 1).Value = 're : change of address - abcd hi test fund the members address has been updated thank you test pensions from : testfund fincorp com sent : 13 october 2020 18 : 14 to : test pensions pensions com subject : change of address - abcd nino - zz123123z dob - 01 01 1400 new address : 10 test road test street aa12 123 tel : + 91 012345678 email abc xyz testmail com could you please update your records with the members new address and contact details could you also previde him with his pension options see attached request kind regards test fund'
$ws.Cells.Item(23, 2).Value = 'Change of Address-ABCD.msg'

$ws.Cells.Item(24, 1).Value = 'abcd xx tst12345678 good morning re xx test pensions ref no tst12345678 sorry this is quite urgent as docs only valid until 01 01 1400 i have tried to contact your office vis phone today still holding after 45 minutes documents received to transfer client benefits can we email across discharge forms signed by client or are originals required no adviser sign off doc was attached is this required can we send across our standard document of can you provide such a document via email asap again is original required in post or can we email i look forward to your response test ifa'
$ws.Cells.Item(24, 2).Value = 'ABCD XX TST12345678.msg'

$ws.Cells.Item(25, 1).Value = 'change of address dear test pension please note that i have moved from : 107 test road test street wa12 123 to : 10 test road test street aa12 123 regards abc + 91 012345678'
$ws.Cells.Item(25, 2).Value = 'ABC_Change of address.msg'

$ws.Cells.Item(26, 1).Value = ' pre - identified user support request reply email address : test ifa fincorp com forename : abc surname : xyz ni number : zz123123z dob : 01 jan 1400 home postcode : aa12 12 employer : fin corp ltd subject : fin corp ltd content : policy 1234567 as i reach retirement age - wish to make arrangements for bringing this pension into payment abc'
$ws.Cells.Item(26, 2).Value = 'FW PRE-IDENTIFIED USER SUPPORT REQUEST 1.msg'

$ws.Cells.Item(27, 1).Value = 'change of address please could you send me a form for change of address my current address is mrs abc 10 test road test street aa12 123 my old address is 09 test road test street ee12 123 regards abc'
$ws.Cells.Item(27, 2).Value = 'Change of addresss.msg'

$ws.Cells.Item(28, 1).Value = 'change of address i have changed my address and now living in xxx my new address is : - 10 test road test street aa12 123 best regards abc'
$ws.Cells.Item(28, 2).Value = 'ABCChange of address.msg'

$ws.Cells.Item(29, 1).Value = 'employee no 132456 hi im an ex employee trying to find all my old pensions from ex employers please could you inform me of steps required in obtaining my contribution while i worked for dsg i am planning on combining all of my pensions so need the data regards abc xyz'
$ws.Cells.Item(29, 2).Value = 'FW EMPLOYEE NO 123456.msg'

$ws.Cells.Item(30, 1).Value = 'fw : change of address from : abc xyz testmail com sent : thursday october 15 2020 3 : 00 : 24 to : test pensions pensions com subject : change of address i would like to update my personal details can someone please contact me as i am not able to contact you via telephone regards abc'
$ws.Cells.Item(30, 2).Value = 'Change of Address  ABC.msg'

$ws.Cells.Item(31, 1).Value = 're : change of personal details name abc further to my previous email my date of birth is 01 january 1400 t i pension reference number tst12345678 national insurance number zz123123z change of address from 70 test road test street wa12 123 new address 10 test road test street aa12 123 mobile + 91 012345678 email abc xyz testmail com<mailto : email jackie barrett123197 gmail com regards abc on wed 14 oct 2020 15 : 01 abc xyz testmail com<mailto : abc xyz testmail com wrote : name abc test group pension scheme member reference number tst12345678 national insurance number zz123123z i wish to advice you on my change of address now; 10 test road test street aa12 123 mobile number + 91 012345678 regards abc - - abc'
$ws.Cells.Item(31, 2).Value = 'Change of personal details.msg'
